# Trials_KFJ.xlsx edit: append 12 new trial rows (374-385) to the "base"
# sheet, extend the AutoFilter/_FilterDatabase/dimension to match, and
# update the saved window state (zoom + selection) on "base" and "info".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("base")

# --- 1. New data rows (374-385), columns A:W (static values) -------------
# Pipe-delimited: A,B,C,D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T,U,V,W
$rowsText = @(
    "300|0|400|0.3|2|0.01|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|0|400|0.6|2|0.01|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|0|400|0.99|2|0.01|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|0|400|0.3|1|0.04|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|0|400|0.6|1|0.04|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|0|400|0.99|1|0.04|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|1|400|0.3|2|0.01|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|1|400|0.6|2|0.01|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|1|400|0.99|2|0.01|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|1|400|0.3|1|0.04|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|1|400|0.6|1|0.04|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1",
    "300|1|400|0.99|1|0.04|0.6|5|7|1|0|30|0|NULL|0|0|NULL|0|0|0|0.2|0.2|T1"
)

$startRow = 374
$nrows = $rowsText.Length
$ncols = 23

$data = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    $parts = $rowsText[$i].Split('|')
    for ($j = 0; $j -lt $ncols; $j++) {
        $s = $parts[$j]
        if ($s -match '^-?[0-9]+(\.[0-9]+)?$') {
            $data[$i, $j] = [double]$s
        }
        else {
            $data[$i, $j] = $s
        }
    }
}

$endRow = $startRow + $nrows - 1
$ws.Range("A$startRow`:W$endRow").Value2 = $data

# --- 2. Formula columns X, Y, Z for the new rows --------------------------
# X: lookup-based trial-dose label
$ws.Range("X$startRow`:X$endRow").Formula = `
    "=CONCATENATE(LOOKUP(D$startRow,info!`$C`$11:`$D`$19), F$startRow*100)"

# Y: sex/replicate label
$ws.Range("Y$startRow`:Y$endRow").Formula = `
    "=IF(AND(B$startRow=0,E$startRow=1),""F1"",IF(AND(B$startRow=0,E$startRow=2),""F2"",IF(AND(B$startRow=1,E$startRow=1),""M1"",IF(AND(B$startRow=1,E$startRow=2),""M2"",""?""))))"

# Z: concatenated label
$ws.Range("Z$startRow`:Z$endRow").Formula = `
    "=CONCATENATE(`$Y$startRow,""-"",`$W$startRow,""-"",`$X$startRow)"

# --- 3. Extend AutoFilter / dimension / _FilterDatabase to the new range -
$ws.AutoFilterMode = $false
$ws.Range("A1:AB$endRow").AutoFilter(1)

foreach ($n in $wb.Names) {
    if ($n.Name -eq "base!_FilterDatabase") {
        $n.RefersTo = "=base!`$A`$1:`$AB`$$endRow"
    }
}

# --- 4. Window state: "info" sheet selection (set before "base" so the ---
#        final active tab stays on "base", matching the target file).
$wsInfo = $wb.Worksheets.Item("info")
$wsInfo.Range("D17").Select()

# --- 5. Window state: "base" sheet zoom + selection -----------------------
$ws.Activate()
$ws.Range("Y391").Select()
$excel.ActiveWindow.Zoom = 90
